$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column (H) is added to the right of the existing "sum" column (G).
# H1 gets the same header formatting as the other header cells (B1:G1): bold
# font, thin box border, centered horizontally and top-aligned vertically.
$xlHAlignCenter = -4108
$xlVAlignTop = -4160
$xlContinuous = 1

$ws.Range("H1").Value = "Save"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = $xlHAlignCenter
$ws.Range("H1").VerticalAlignment = $xlVAlignTop
$ws.Range("H1").Borders.LineStyle = $xlContinuous

# Fill in the "Save" column values for the data rows (2-6)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
